$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1) Insert two fresh rows at the top of the data block (rows 2 & 3) ---
# Shifts the old rows 2-5 down to rows 4-7, carrying their formatting along.
$ws.Range("A2:A3").EntireRow.Insert()

# --- 2) New row 2: マッチングアプリのAIレコメンド構築 ---
$ws.Range("A2").Value2 = "2025-10-09 18:24:13"
$ws.Range("B2").Value2 = "マッチングアプリのAIレコメンド構築"
$ws.Range("C2").Value2 = "システム開発"
$ws.Range("D2").Value2 = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E2").Value2 = "期限情報なし"
$ws.Range("F2").Value2 = "https://www.lancers.jp/work/detail/5410515"
$ws.Range("G2").Value2 = 338
$ws.Range("H2").Value2 = "🔥AI,Ai ◇アプリ"

# --- 3) New row 3: 【バックエンド開発】ポータルサイトの予約情報管理システム構築 ---
$ws.Range("A3").Value2 = "2025-10-09 18:24:13"
$ws.Range("B3").Value2 = "【バックエンド開発】ポータルサイトの予約情報管理システム構築"
$ws.Range("C3").Value2 = "システム開発"
$ws.Range("D3").Value2 = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E3").Value2 = "期限情報なし"
$ws.Range("F3").Value2 = "https://www.lancers.jp/work/detail/5410302"
$ws.Range("G3").Value2 = 143
$ws.Range("H3").Value2 = "◆開発 ◇サイト"

# --- 4) Refresh the "seen at" timestamp on the previously-known rows (now 4-7) ---
$ws.Range("A4").Value2 = "2025-10-09 18:24:13"
$ws.Range("A5").Value2 = "2025-10-09 18:24:13"
$ws.Range("A6").Value2 = "2025-10-09 18:24:13"
$ws.Range("A7").Value2 = "2025-10-09 18:24:13"

# --- 5) Rebuild the hyperlinks for column F so every URL row (2-7) carries
#        exactly one correct hyperlink (row-insert does not renumber the
#        sheet's existing <hyperlinks> refs, so drop them all and re-add). ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5410515")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5410302")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5251319")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5409967")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5410017")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5410127")
